# js to validate form
$wb = $excel.ActiveWorkbook

# --- Rename sheet "6.5" to "6.5+8.5" ---
$wsA = $wb.Worksheets.Item("24.4")
$wsB = $wb.Worksheets.Item("6.5")
$wsB.Name = "6.5+8.5"

# --- Fill in the two new rows of content on the "6.5+8.5" sheet ---
$wsB.Range("C7").Value = "CSS - HTML: Cách nắm đầu các thẻ đẻ chỉnh css"
$wsB.Range("C8").Value = "JS - hTML: Cách nắm đầu các thẻ theo object để chỉnh sửa theo js"

# --- Update scroll position / selection on "6.5+8.5" ---
$wsB.Activate()
[void]$wsB.Application.Goto($wsB.Range("A13"), $true)
[void]$wsB.Range("C11").Select()

# --- Switch the active sheet back to "24.4" (clears activeTab / tabSelected on the other sheet) ---
$wsA.Activate()
[void]$wsA.Range("C28").Select()
